$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.572.27'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.57%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.640.68'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.86%  '

$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '212.69'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.41%  '

$ws.Range("E6").Value = '  +4.69%  '

$ws.Range("E7").Value = '  -0.03%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '23.08'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -4.00%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.257'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.75%  '

$ws.Range("E10").Value = '  -0.62%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0890'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.38%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.871.48'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.92%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.639.52'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.04%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.03'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.12%  '

$ws.Range("E15").Value = '  -2.02%  '

$ws.Range("E16").Value = '  -2.72%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '27.573.81'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.54%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '229.56'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.03%  '

$ws.Range("B19").Value = 'ShibaInu'
$ws.Range("C19").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.0₃0724'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.24%  '

$ws.Range("B20").Value = 'Chainlink'
$ws.Range("C20").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.64'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.85%  '

$ws.Range("E21").Value = '  +0.06%  '

$ws.Range("E22").Value = '  -1.71%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.94'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +6.74%  '

$ws.Range("E24").Value = '  -3.53%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '149.35'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.76%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.96'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -3.53%  '

$ws.Range("E28").Value = '  -0.15%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.62'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.97%  '

$ws.Range("E30").Value = '  -1.07%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0486'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.14%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.30'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.16%  '

$ws.Range("E33").Value = '  +1.03%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.426.81'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.65%  '

$ws.Range("E35").Value = '  +2.05%  '

$ws.Range("E36").Value = '  -2.07%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.573'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.28%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.879'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.45%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0167'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.57%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.903'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +15.55%  '

$ws.Range("E41").Value = '  -1.80%  '

$ws.Range("E42").Value = '  +0.02%  '

$ws.Range("E43").Value = '  -1.13%  '

$ws.Range("E44").Value = '  +1.03%  '

$ws.Range("E45").Value = '  +1.57%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '65.08'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.45%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.781.53'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.87%  '

$ws.Range("E48").Value = '  -2.07%  '

$ws.Range("E49").Value = '  -2.32%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0₆0107'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.76%  '

$ws.Range("E51").Value = '  -2.54%  '
